# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on the active sheet to reflect the regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 3
    11 = 2
    12 = 2
    13 = 0
    14 = 3
    15 = 3
    16 = 0
    17 = 1
    18 = 1
    19 = 0
    20 = 1
    21 = 3
    22 = 1
    23 = 2
    24 = 2
    26 = 2
    27 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
